# Bump the "想去人数" (want-to-go count) column F for 5 events by +1,
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": F3, F7, F12, F15, F18
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 679
$ws1.Range("F7").Value  = 10984
$ws1.Range("F12").Value = 10829
$ws1.Range("F15").Value = 24
$ws1.Range("F18").Value = 82

# Sheet "全部类型": F3, F10, F15, F18, F21 (same five events, different rows)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 679
$ws4.Range("F10").Value = 10984
$ws4.Range("F15").Value = 10829
$ws4.Range("F18").Value = 24
$ws4.Range("F21").Value = 82
